$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the note text in B9 to reflect the newly calculated total, and
# align it like the rest of the column (center, matching the other day
# rows) instead of the old left-aligned "note" style.
$ws.Range("B9").Value = "13 Hours 25 Minutes"
$ws.Range("B9").HorizontalAlignment = -4108  # xlCenter

# Move the active selection to C12 (where the author left off).
$ws.Range("C12").Select()
